$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range first (rows 3-9 need to go away entirely)
$ws.Range("A1:B9").Clear()

# New header row
$ws.Range("A1").Value = "curriculum"
$ws.Range("B1").Value = "english_name"
$ws.Range("C1").Value = "placement_id"
$ws.Range("D1").Value = "placement_year"
$ws.Range("E1").Value = "company_name"

# New data row
$ws.Range("A2").Value = "BBA(IS)"
$ws.Range("B2").Value = "Bob Doe"
$ws.Range("C2").Value = 4
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2022"
